# DICCIONARIO-DE-DATOS.xlsx — "Add files via upload" revision
#
# The original author re-saved the workbook after renaming several field
# names in the PRESTAMO / USUARIOS entity tables (rows 44-48 and 55-58 on
# Hoja1) to a consistent CamelCase / underscore convention, and the
# selection/scroll position moved as a side effect of editing near the
# bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- PRESTAMO entity (rows 44-48) : rename CAMPO column values ---------
$ws.Range("C44").Value = "IdPrestamo"        # was: idprestamo
$ws.Range("C45").Value = "Fecha_Prestamo"    # was: fecha_prestamo
$ws.Range("C46").Value = "Fecha_Entrega"     # was: fecha_entrega
$ws.Range("C47").Value = "Ejemplares_sol"    # was: cantidad_ejemplares_solicitados
$ws.Range("C48").Value = "Desc_Prestamo"     # was: descripción_prestamo

# Row 47 no longer needs its taller custom height once the label got
# shorter — reset it back to the sheet's default row height.
$ws.Rows.Item(47).AutoFit()

# --- USUARIOS entity (rows 55-58) : rename CAMPO column values ---------
$ws.Range("C55").Value = "IdUsuario"  # was: idusuario
$ws.Range("C56").Value = "Correo"     # was: correo
$ws.Range("C57").Value = "Usuario"    # was: usuario
$ws.Range("C58").Value = "Contraseña" # was: contraseña

# --- Cursor ended up on D58 after the last edit -------------------------
$null = $ws.Range("D58").Select()
